$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix typo: "hope page" -> "home page"
$ws.Range("B3").Value = "ketika sedang meingisi dan belum selesai, kemudia user klik tombol summary, maka harus redirect to home page applicant, bukan halaman error"

Write-Host "done"
